$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row data: RowNum, Coin, Link, Price, Volume(1h)
$data = @(
    ,@(2, "Bitcoin", "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc", "67.768.20", "  +0.04%  ")
    ,@(3, "Ethereum", "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth", "3.808.82", "  +0.49%  ")
    ,@(4, "TetherUSD", "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt", "1.00", "  +0.01%  ")
    ,@(5, "BNB", "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb", "603.29", "  +1.23%  ")
    ,@(6, "Solana", "https://coinranking.com/coin/zNZHO_Sjf+solana-sol", "166.04", "  -0.67%  ")
    ,@(7, "USDC", "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc", "1.00", "  +0.03%  ")
    ,@(8, "XRP", "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp", "0.518", "  -0.24%  ")
    ,@(9, "Dogecoin", "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge", "0.159", "  -0.03%  ")
    ,@(10, "Cardano", "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada", "0.453", "  +0.68%  ")
    ,@(11, "Toncoin", "https://coinranking.com/coin/67YlI0K1b+toncoin-ton", "6.35", "  +0.83%  ")
    ,@(12, "ShibaInu", "https://coinranking.com/coin/xz24e0BjL+shibainu-shib", "0.0000250", "  -1.40%  ")
    ,@(13, "Avalanche", "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax", "35.95", "  -0.27%  ")
    ,@(14, "WrappedliquidstakedEther2.0", "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth", "4.451.43", "  +0.70%  ")
    ,@(15, "WrappedEther", "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth", "3.795.61", "  +0.85%  ")
    ,@(16, "WrappedBTC", "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc", "67.805.76", "  +0.15%  ")
    ,@(17, "Chainlink", "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link", "18.36", "  -0.81%  ")
    ,@(18, "TRON", "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx", "0.113", "  +1.81%  ")
    ,@(19, "Polkadot", "https://coinranking.com/coin/25W7FG7om+polkadot-dot", "7.07", "  +0.45%  ")
    ,@(20, "BitcoinCash", "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch", "464.08", "  +0.79%  ")
    ,@(21, "Uniswap", "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni", "9.82", "  -2.01%  ")
    ,@(22, "Polygon", "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic", "0.701", "  +0.61%  ")
    ,@(23, "Litecoin", "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc", "83.32", "  -0.06%  ")
    ,@(24, "PEPE", "https://coinranking.com/coin/03WI8NQPF+pepe-pepe", "0.0000146", "  -4.83%  ")
    ,@(25, "InternetComputer(DFINITY)", "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp", "12.12", "  +0.69%  ")
    ,@(26, "Fetch.AI", "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet", "2.12", "  -0.74%  ")
    ,@(27, "RenderToken", "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr", "10.01", "  -0.19%  ")
    ,@(28, "Dai", "https://coinranking.com/coin/MoTuySvg7+dai-dai", "1.00", "  -0.13%  ")
    ,@(29, "WrappedeETH", "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth", "3.961.00", "  +0.75%  ")
    ,@(30, "PancakeSwap", "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake", "2.78", "  +0.22%  ")
    ,@(31, "NEARProtocol", "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near", "7.42", "  +2.71%  ")
    ,@(32, "ImmutableX", "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx", "2.21", "  -1.30%  ")
    ,@(33, "EthereumClassic", "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc", "29.43", "  -0.86%  ")
    ,@(34, "Binance-PegBSC-USD", "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd", "1.00", "  +0.10%  ")
    ,@(35, "Aptos", "https://coinranking.com/coin/HGYj5JCv5+aptos-apt", "9.06", "  -0.56%  ")
    ,@(36, "Hedera", "https://coinranking.com/coin/jad286TjB+hedera-hbar", "0.0997", "  -0.52%  ")
    ,@(37, "Kaspa", "https://coinranking.com/coin/V8GxkwWow+kaspa-kas", "0.138", "  +0.45%  ")
    ,@(38, "Mantle", "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt", "0.996", "  +0.04%  ")
    ,@(39, "Filecoin", "https://coinranking.com/coin/ymQub4fuB+filecoin-fil", "5.81", "  +0.71%  ")
    ,@(40, "dogwifhat", "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif", "3.23", "  -3.42%  ")
    ,@(41, "FirstDigitalUSD", "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd", "1.00", "  +0.07%  ")
    ,@(42, "USDe", "https://coinranking.com/coin/exbfr2U-0+usde-usde", "1.00", "  +0.01%  ")
    ,@(43, "Arweave", "https://coinranking.com/coin/7XWg41D1+arweave-ar", "44.58", "  -3.42%  ")
    ,@(44, "OKB", "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb", "47.74", "  -0.85%  ")
    ,@(45, "TheGraph", "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt", "0.299", "  -0.09%  ")
    ,@(46, "EnergySwap", "https://coinranking.com/coin/SbWqqTui-+energyswap-ens", "28.02", "  +4.89%  ")
    ,@(47, "Monero", "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr", "151.67", "  +1.61%  ")
    ,@(48, "ONDO", "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo", "1.39", "  +11.80%  ")
    ,@(49, "Cosmos", "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom", "8.34", "  +0.12%  ")
    ,@(50, "Stacks", "https://coinranking.com/coin/mMPrMcB7+stacks-stx", "1.84", "  +1.19%  ")
    ,@(51, "Bittensor", "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao", "389.78", "  -1.19%  ")
)

foreach ($item in $data) {
    $row = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    # Force Price column to remain plain text (avoid numeric auto-conversion of values like "1.00" or "67.768.20")
    $ws.Cells.Item($row, 4).NumberFormat = "@"
    $ws.Cells.Item($row, 4).Value = $item[3]
    $ws.Cells.Item($row, 5).Value = $item[4]
}

$wb.Save()